$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks first (their underlying text stays intact); we will
# recreate them afterwards once the new row has been inserted and everything
# has shifted into its final position.
$ws.Hyperlinks.Delete()

# Insert a new blank row above row 5; this pushes the existing project rows
# (old rows 5-8) down to rows 6-9.
$ws.Rows("5:5").Insert()

# Populate the new row with the "Packet Sniffer" project.
$ws.Range("A5").Value = "Packet Sniffer"
$ws.Range("C5").Value = "Python"
$ws.Range("D5").Value = "https://github.com/c-l-scholl/python-packet-sniffer"
$ws.Range("B5").Value = "Used Python sockets and structs to create a basic packet sniffer. Based on a tutorial by thenewboston"

# Give the new row the same cell formatting used by the other project rows
# (row 6, directly below, already carries the correct per-cell styles).
$ws.Range("A6:D6").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row heights seen in the final layout.
$ws.Rows("5:5").RowHeight = 29.4
$ws.Rows("6:6").RowHeight = 43.2
$ws.Rows("7:7").RowHeight = 43.2
$ws.Rows("8:8").RowHeight = 28.8
$ws.Rows("9:9").RowHeight = 28.8

# Recreate the github-link hyperlinks, now that everything is in its final
# row position. The original four keep their original relative order (so
# they reuse the same relationship-id ordering as before), and the brand
# new Packet Sniffer link is appended last.
$ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/c-l-scholl/stp-vue-fb")
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/mac-comp128-s22/128-project-ben-and-camden2")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/c-l-scholl/YT-home-page-UI-replica")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://github.com/c-l-scholl/discord-bot")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/c-l-scholl/python-packet-sniffer")

# Adding a hyperlink re-applies its own "Hyperlink" style record; reapply the
# pre-existing per-column formatting on top so the cells keep reusing the
# workbook's original style entries instead of duplicating new ones.
$ws.Range("D9").Copy()
$ws.Range("D5:D9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the final selection shown in the workbook.
$ws.Range("B5").Select()
